$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.00000000000000000000000000502558222194548773272546
$ws.Range("E2").Value = 0.00000000000000000000000000502558222194548773272546

$ws.Range("D3").Value = 0.99999999999983657517077517695724964141845703125000
$ws.Range("E3").Value = 0.99999999999983657517077517695724964141845703125000

$ws.Range("D4").Value = 0.00000000002154099554096399072113612055479531020782
$ws.Range("E4").Value = 0.99999999997845900878701286274008452892303466796875

$ws.Range("D5").Value = 0.99989724960335013559387107306974940001964569091797
$ws.Range("E5").Value = 0.00010275039664986440612892693025059998035430908203

$ws.Range("D6").Value = 0.00046451617666640188097915142861893400549888610840
$ws.Range("E6").Value = 0.99953548382333357036344523294246755540370941162109

$ws.Range("D7").Value = 0.00000000000000063100168532281528365573511855720026
$ws.Range("E7").Value = 0.99999999999999933386618522490607574582099914550781

$ws.Range("D8").Value = 0.00000000000006608936468993912348680136611569946623
$ws.Range("E8").Value = 0.99999999999993394173003480318584479391574859619141

$ws.Range("F8").Value = 18.14642906188964843750000000000000000000000000000000
